$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.07 = 24162.74 pesos`n✅ 24162.74 pesos = 6.04 = 969.86 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 164.8
$wsTasas.Range("O10").Value = 3982.02
$wsTasas.Range("N12").Value = 4000
$wsTasas.Range("O12").Value = 160.555
